# Reorder the two worksheets so "review_info" becomes the first tab and
# "hotel_info" becomes the second tab (previously hotel_info was first).
$wb = $excel.ActiveWorkbook

$wsReview = $wb.Worksheets.Item("review_info")
$wsReview.Move($wb.Worksheets.Item(1))

# Worksheet object references resolved before the Move can become stale,
# so re-fetch "hotel_info" by name now that the tab order has changed.
$wsHotel = $wb.Worksheets.Item("hotel_info")

# Insert a new "State" column into hotel_info, right after "Hotel_Name"
# and before "City", shifting the existing columns (City, Zip, ...) right.
$wsHotel.Range("C1").EntireColumn.Insert()
$wsHotel.Range("C1").Value = "State"
$wsHotel.Range("C2").Value = "Louisiana"
